$wb = $excel.ActiveWorkbook

$targets = @("展览", "全部类型")
foreach ($name in $targets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2123
    $ws.Range("F3").Value = 616
    $ws.Range("F4").Value = 1506
    $ws.Range("F5").Value = 7161
    $ws.Range("F7").Value = 143
}
